$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2026-02-26 14:28:08"

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 10).Value = $newTimestamp
}
